# Commit: "Added further education, and presentations for VRI.
# Could not recover footer"
#
# For this workbook (award.xlsx / sheet "education"), the only
# semantic change is that the funding body for the "Francisco José de
# Caldas Scholarship for Doctoral Studies" award (row 7, column C) was
# renamed from its old name "Colciencias" to its current name
# "Minciencias" (Colombia's science ministry was renamed from
# Colciencias to Minciencias in 2019).
#
# Updating that single cell is enough to also reproduce the
# sharedStrings.xml cleanup seen in the diff: "Colciencias" drops out
# (no longer referenced), and the already-unused "Best overall
# performance in the MSc." (trailing-period variant) string is no
# longer carried along, while the new "Minciencias" string is appended
# at the end of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("education")

$ws.Range("C7").Value = "Minciencias"

# The author's last selection in the sheet moved to C15 before saving.
$ws.Range("C15").Select()
